$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 7, 8, 9 with new values ---

# Row 7 (Q5, N=11 -> 9)
$ws.Range("B7").Value = -0.2848744657242845
$ws.Range("C7").Value = 1.081999780872457
$ws.Range("D7").Value = 1.826018359938472
$ws.Range("E7").Value = 1.351302467968764
$ws.Range("F7").Value = 1.401061387326184
$ws.Range("G7").Value = 9

# Row 8 (Q6, N=12 -> 6)
$ws.Range("B8").Value = -0.660147827454144
$ws.Range("C8").Value = 0.8867455175475859
$ws.Range("D8").Value = 0.9350091296531003
$ws.Range("E8").Value = 0.9669587011103941
$ws.Range("F8").Value = 0.7739875778543277
$ws.Range("G8").Value = 6

# Row 9 (Q7, N=13 -> 3), F9 newly populated
$ws.Range("B9").Value = -0.01627907676619128
$ws.Range("C9").Value = 0.5255089814372506
$ws.Range("D9").Value = 0.3174133864681257
$ws.Range("E9").Value = 0.5633945211555804
$ws.Range("F9").Value = 0.689726443738131
$ws.Range("G9").Value = 3

# --- Add new row 10 (Q8) ---
# Copy formatting from A9 (bold font, thin border, centered) onto the new A10 cell
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = "Q8"

$ws.Range("B10").Value = -0.2577465226711695
$ws.Range("C10").Value = 0.2577465226711695
$ws.Range("D10").Value = 0.06643326994907969
$ws.Range("E10").Value = 0.2577465226711695
$ws.Range("G10").Value = 1
